$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J (row 1), matching the formatting
# already used by the other header cells (e.g. H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data values for row 2 in columns I and J.
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 8
